$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: the order in which brand-new requirement texts are first written
# matters for shared-string ordering, so B6..B15/B2..B5 are touched in the
# same sequence the authoring tool produced them in.

$ws.Range("B6").Value = "Protótipo possuir nome da empresa/grupo e a nossa proposta"
$ws.Range("C6").Value = "essencial"

$ws.Range("B7").Value = "Protótipo explicar o nosso projeto principal"
$ws.Range("B7").Font.ThemeFont = 1
$ws.Range("C7").Value = "essencial"

$ws.Range("B8").Value = "Protótipo possuir uma área ""sobre nós"""
$ws.Range("B8").Font.ThemeFont = 1
$ws.Range("C8").Value = "importante"

$ws.Range("B9").Value = "Protótipo possuir uma área ""curiosidades sobre o mundo orgânico"""
$ws.Range("B9").Font.ThemeFont = 1
$ws.Range("C9").Value = "importante"

$ws.Range("B10").Value = "Protótipo possuir uma área ""e-commerce"""
$ws.Range("B10").Font.ThemeFont = 1
$ws.Range("C10").Value = "desejavel"

$ws.Range("B11").Value = "Botão que direciona à homepage do protótipo"
$ws.Range("B11").Font.ThemeFont = 1
$ws.Range("C11").Value = "essencial"

$ws.Range("B12").Value = "Protótipo mostrar o simulador financeiro"
$ws.Range("B12").Font.ThemeFont = 1
$ws.Range("C12").Value = "essencial"

$ws.Range("B2").Value = "Protótipo com tela de cadastro para o usuário"

$ws.Range("B3").Value = "Protótipo com tela  de login para o usuário"

$ws.Range("B4").Value = "Usuário ter acesso aos gráficos forncecidos pelo Arduíno através do protótipo"

$ws.Range("B13").Value = "O usuário pode ver dados do banco de dados das plantas através do protótipo"
$ws.Range("B13").Font.ThemeFont = 1
$ws.Range("C13").Value = "importante"

$ws.Range("B5").Value = "Protótipo possuir a Logo da equipe em todas as telas"
$ws.Range("B5").Font.ThemeFont = 1

$ws.Range("B14").Value = "Arduíno conseguir controlar a estufa"
$ws.Range("B14").Font.ThemeFont = 1
$ws.Range("C14").Value = "essencial"

$ws.Range("B15").Value = "Site completo no futuro"
$ws.Range("B15").Font.ThemeFont = 1
$ws.Range("C15").Value = "essencial"

# --- Column G (second table) text refreshed (content identical, just kept in sync) ---
$ws.Range("G3").Value = "Banco de dados para cadastro"
$ws.Range("G4").Value = "Ter conexão com a internet"
$ws.Range("G5").Value = "Banco de dados conectado à nuvem"
$ws.Range("G6").Value = "Sensores conversarem com o arduíno"
$ws.Range("G7").Value = "Arduíno conversar com o sistema"

# --- Column B width widened to fit the longer requirement text ---
$ws.Columns.Item(2).ColumnWidth = 69

# --- Selection moved to B19 ---
$ws.Range("B19").Select()
